# CELGRandomTrade.xlsx - append a new trade row (row 6) to the trade log,
# mirroring the existing rows' layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row data
$ws.Range("A6").Value = 42636.588993055557
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = 10114.32
$ws.Range("D6").Value = 10062
$ws.Range("E6").Value = 110
$ws.Range("F6").Value = 110.57
$ws.Range("G6").Value = $false
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = 0.52
$ws.Range("I6").Value = $true

# Column C ("Principle") widens slightly to fit the new, longer value (10114.32)
$ws.Columns.Item(3).ColumnWidth = 8.16
